# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# to the latest pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 8).Value = 20987.25
$ws.Cells.Item(21, 9).Value = 20987.25
$ws.Cells.Item(21, 11).Value = 20987.25
$ws.Cells.Item(21, 13).Value = -20519.25

$ws.Cells.Item(23, 8).Value = 20987.25
$ws.Cells.Item(23, 9).Value = 20987.25
$ws.Cells.Item(23, 11).Value = 20987.25
$ws.Cells.Item(23, 13).Value = -20753.25

$ws.Cells.Item(74, 8).Value = 14306.772
$ws.Cells.Item(74, 9).Value = 17092.8
$ws.Cells.Item(74, 10).Value = 8336.714
$ws.Cells.Item(74, 11).Value = 17092.8
$ws.Cells.Item(74, 12).Value = 8336.714
$ws.Cells.Item(74, 13).Value = -16156.8
$ws.Cells.Item(74, 14).Value = -10208.714

$ws.Cells.Item(77, 8).Value = 14306.772
$ws.Cells.Item(77, 9).Value = 17092.8
$ws.Cells.Item(77, 10).Value = 8336.714
$ws.Cells.Item(77, 11).Value = 85464
$ws.Cells.Item(77, 12).Value = 41683.57
$ws.Cells.Item(77, 13).Value = -80784
$ws.Cells.Item(77, 14).Value = -51043.57

$ws.Cells.Item(98, 8).Value = 2530.3428
$ws.Cells.Item(98, 9).Value = 2488.2727
$ws.Cells.Item(98, 11).Value = 2488.2727
$ws.Cells.Item(98, 13).Value = -990.2727

$ws.Cells.Item(122, 8).Value = 2530.3428
$ws.Cells.Item(122, 9).Value = 2488.2727
$ws.Cells.Item(122, 11).Value = 7464.8181
$ws.Cells.Item(122, 13).Value = -5014.8181

$ws.Cells.Item(132, 8).Value = 2619479
$ws.Cells.Item(132, 9).Value = 5397.4062
$ws.Cells.Item(132, 10).Value = 44444788
$ws.Cells.Item(132, 11).Value = 16192.2186
$ws.Cells.Item(132, 12).Value = 133334364
$ws.Cells.Item(132, 13).Value = -13662.2186
$ws.Cells.Item(132, 14).Value = -133339424

$ws.Cells.Item(137, 8).Value = 1253033.8
$ws.Cells.Item(137, 9).Value = 1726986.5
$ws.Cells.Item(137, 11).Value = 5180959.5
$ws.Cells.Item(137, 13).Value = -5178409.5

$ws.Cells.Item(138, 8).Value = 2312.889
$ws.Cells.Item(138, 9).Value = 1219.5
$ws.Cells.Item(138, 10).Value = 4499.6665
$ws.Cells.Item(138, 11).Value = 3658.5
$ws.Cells.Item(138, 12).Value = 13498.9995
$ws.Cells.Item(138, 13).Value = 1481.5
$ws.Cells.Item(138, 14).Value = -23778.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2614.743
$ws.Cells.Item(32, 9).Value = 2632.4412
$ws.Cells.Item(32, 10).Value = 2013
$ws.Cells.Item(32, 11).Value = 2632.4412
$ws.Cells.Item(32, 12).Value = 2013
$ws.Cells.Item(32, 13).Value = -2345.4412
$ws.Cells.Item(32, 14).Value = -2587

$ws.Cells.Item(61, 8).Value = 2985.3684
$ws.Cells.Item(61, 9).Value = 2518.5833
$ws.Cells.Item(61, 11).Value = 2518.5833
$ws.Cells.Item(61, 13).Value = -2306.5833

$ws.Cells.Item(102, 8).Value = 3949.5833
$ws.Cells.Item(102, 10).Value = 4468.375
$ws.Cells.Item(102, 12).Value = 4468.375
$ws.Cells.Item(102, 14).Value = -7712.375

$ws.Cells.Item(122, 8).Value = 8077.375
$ws.Cells.Item(122, 9).Value = 8929.368
$ws.Cells.Item(122, 11).Value = 26788.104
$ws.Cells.Item(122, 13).Value = -24338.104

$ws.Cells.Item(136, 8).Value = 2985.3684
$ws.Cells.Item(136, 9).Value = 2518.5833
$ws.Cells.Item(136, 11).Value = 7555.749899999999
$ws.Cells.Item(136, 13).Value = -5005.749899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 169533.33
$ws.Cells.Item(7, 9).Value = 551.5
$ws.Cells.Item(7, 11).Value = 551.5
$ws.Cells.Item(7, 13).Value = -438.5

$ws.Cells.Item(75, 8).Value = 14623.571
$ws.Cells.Item(75, 9).Value = 15060.833
$ws.Cells.Item(75, 11).Value = 15060.833
$ws.Cells.Item(75, 13).Value = -14124.833

$ws.Cells.Item(78, 8).Value = 14623.571
$ws.Cells.Item(78, 9).Value = 15060.833
$ws.Cells.Item(78, 11).Value = 45182.499
$ws.Cells.Item(78, 13).Value = -40502.499

$ws.Cells.Item(94, 8).Value = 153848770
$ws.Cells.Item(94, 9).Value = 333334900
$ws.Cells.Item(94, 10).Value = 3497.4285
$ws.Cells.Item(94, 11).Value = 333334900
$ws.Cells.Item(94, 12).Value = 3497.4285
$ws.Cells.Item(94, 13).Value = -333334449
$ws.Cells.Item(94, 14).Value = -4399.4285

$ws.Cells.Item(99, 8).Value = 3813.7693
$ws.Cells.Item(99, 9).Value = 2676.4443
$ws.Cells.Item(99, 11).Value = 2676.4443
$ws.Cells.Item(99, 13).Value = -1178.4443

$ws.Cells.Item(134, 8).Value = 3262.3447
$ws.Cells.Item(134, 9).Value = 3072.5
$ws.Cells.Item(134, 11).Value = 9217.5
$ws.Cells.Item(134, 13).Value = -6682.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2556932.8
$ws.Cells.Item(31, 9).Value = 3919.4333
$ws.Cells.Item(31, 10).Value = 6588006.5
$ws.Cells.Item(31, 11).Value = 3919.4333
$ws.Cells.Item(31, 12).Value = 6588006.5
$ws.Cells.Item(31, 13).Value = -3624.4333
$ws.Cells.Item(31, 14).Value = -6588596.5

$ws.Cells.Item(34, 8).Value = 2556932.8
$ws.Cells.Item(34, 9).Value = 3919.4333
$ws.Cells.Item(34, 10).Value = 6588006.5
$ws.Cells.Item(34, 11).Value = 3919.4333
$ws.Cells.Item(34, 12).Value = 6588006.5
$ws.Cells.Item(34, 13).Value = -3717.4333
$ws.Cells.Item(34, 14).Value = -6588410.5

$ws.Cells.Item(107, 8).Value = 2174624.5
$ws.Cells.Item(107, 9).Value = 3333828.8
$ws.Cells.Item(107, 10).Value = 1116.625
$ws.Cells.Item(107, 11).Value = 3333828.8
$ws.Cells.Item(107, 12).Value = 1116.625
$ws.Cells.Item(107, 13).Value = -3331908.8
$ws.Cells.Item(107, 14).Value = -4956.625

$ws.Cells.Item(133, 8).Value = 79599.3
$ws.Cells.Item(133, 10).Value = 79599.3
$ws.Cells.Item(133, 12).Value = 79599.3
$ws.Cells.Item(133, 14).Value = -84659.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 1576.75
$ws.Cells.Item(2, 10).Value = 2113
$ws.Cells.Item(2, 12).Value = 12678
$ws.Cells.Item(2, 14).Value = -12904

$ws.Cells.Item(4, 8).Value = 28294786
$ws.Cells.Item(4, 9).Value = 30035460
$ws.Cells.Item(4, 11).Value = 90106380
$ws.Cells.Item(4, 13).Value = -90106268

$ws.Cells.Item(131, 8).Value = 11437.5
$ws.Cells.Item(131, 9).Value = 15420.471
$ws.Cells.Item(131, 10).Value = 1764.5714
$ws.Cells.Item(131, 11).Value = 46261.413
$ws.Cells.Item(131, 12).Value = 5293.7142
$ws.Cells.Item(131, 13).Value = -41221.413
$ws.Cells.Item(131, 14).Value = -15373.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(31, 8).Value = 653
$ws.Cells.Item(31, 9).Value = 653
$ws.Cells.Item(31, 11).Value = 653
$ws.Cells.Item(31, 13).Value = -361

$ws.Cells.Item(37, 8).Value = 653
$ws.Cells.Item(37, 9).Value = 653
$ws.Cells.Item(37, 11).Value = 653
$ws.Cells.Item(37, 13).Value = -376

$ws.Cells.Item(70, 8).Value = 18603528
$ws.Cells.Item(70, 9).Value = 27895578
$ws.Cells.Item(70, 11).Value = 27895578
$ws.Cells.Item(70, 13).Value = -27895308

$ws.Cells.Item(73, 8).Value = 18603528
$ws.Cells.Item(73, 9).Value = 27895578
$ws.Cells.Item(73, 11).Value = 27895578
$ws.Cells.Item(73, 13).Value = -27894642

$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 12).Value = 0
$ws.Cells.Item(101, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 1988.04
$ws.Cells.Item(132, 9).Value = 1994.2727
$ws.Cells.Item(132, 11).Value = 5982.8181
$ws.Cells.Item(132, 13).Value = -3452.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 3394
$ws.Cells.Item(93, 9).Value = 2858.6667
$ws.Cells.Item(93, 10).Value = 5000
$ws.Cells.Item(93, 11).Value = 2858.6667
$ws.Cells.Item(93, 12).Value = 5000
$ws.Cells.Item(93, 13).Value = -1610.6667
$ws.Cells.Item(93, 14).Value = -7496

$ws.Cells.Item(96, 8).Value = 49990
$ws.Cells.Item(96, 10).Value = 49990
$ws.Cells.Item(96, 12).Value = 49990
$ws.Cells.Item(96, 14).Value = -55482

$ws.Cells.Item(122, 8).Value = 5011.231
$ws.Cells.Item(122, 9).Value = 5799.5713
$ws.Cells.Item(122, 10).Value = 4091.5
$ws.Cells.Item(122, 11).Value = 17398.7139
$ws.Cells.Item(122, 12).Value = 12274.5
$ws.Cells.Item(122, 13).Value = -14948.7139
$ws.Cells.Item(122, 14).Value = -17174.5

$ws.Cells.Item(135, 8).Value = 103999.75
$ws.Cells.Item(135, 10).Value = 103999.75
$ws.Cells.Item(135, 12).Value = 103999.75
$ws.Cells.Item(135, 14).Value = -114139.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(21, 8).Value = 5000
$ws.Cells.Item(21, 10).Value = 5000
$ws.Cells.Item(21, 12).Value = 5000
$ws.Cells.Item(21, 14).Value = -5470

$ws.Cells.Item(35, 8).Value = 5000
$ws.Cells.Item(35, 10).Value = 5000
$ws.Cells.Item(35, 12).Value = 5000
$ws.Cells.Item(35, 14).Value = -5580

$ws.Cells.Item(113, 8).Value = 529.23254
$ws.Cells.Item(113, 9).Value = 451.3871
$ws.Cells.Item(113, 10).Value = 730.3333
$ws.Cells.Item(113, 11).Value = 1354.1613
$ws.Cells.Item(113, 12).Value = 2190.9999
$ws.Cells.Item(113, 13).Value = 815.8387
$ws.Cells.Item(113, 14).Value = -6530.9999

$ws.Cells.Item(122, 8).Value = 11906788
$ws.Cells.Item(122, 9).Value = 2186.0527
$ws.Cells.Item(122, 11).Value = 6558.158100000001
$ws.Cells.Item(122, 13).Value = -4108.158100000001

$ws.Cells.Item(132, 8).Value = 4809.6875
$ws.Cells.Item(132, 9).Value = 4551
$ws.Cells.Item(132, 11).Value = 13653
$ws.Cells.Item(132, 13).Value = -11123

$ws.Cells.Item(136, 8).Value = 14488.318
$ws.Cells.Item(136, 9).Value = 15737.974
$ws.Cells.Item(136, 11).Value = 47213.922
$ws.Cells.Item(136, 13).Value = -44663.922
